$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.784.35'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.515.11'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.07%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.73'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +0.90%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.513.43'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('E10').Value = '  +1.16%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.167'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.47%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.362'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.52%  '
$ws.Range('E13').Value = '  +1.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.980.71'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '69.584.06'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.25%  '
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.513.57'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.35'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.64'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '352.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('E22').Value = '  -0.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.99'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '71.04'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.75%  '
$ws.Range('E26').Value = '  -1.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.87'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.48%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.995'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0896'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.26%  '
$ws.Range('E31').Value = '  +0.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '463.81'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.23'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.98%  '
$ws.Range('E34').Value = '  -1.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '159.38'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.02%  '
$ws.Range('E37').Value = '  +0.77%  '
$ws.Range('E38').Value = '  +1.03%  '
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.70'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.46%  '
$ws.Range('E43').Value = '  -1.50%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '38.28'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.23'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.50%  '
$ws.Range('B46').Value = 'ImmutableX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.11'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '143.06'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.49'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.523'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0735'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.89%  '
$ws.Range('E51').Value = '  +3.75%  '
